$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Column C raw value corrections (rows 12-14) ---
$ws.Range("C12").Value = 1427602155.6800008
$ws.Range("D12").Value = 1398800851

$ws.Range("C13").Value = 340926355.51000023
$ws.Range("D13").Value = 338965688.89999998

$ws.Range("C14").Value = 985666.54999992996
$ws.Range("D14").Value = 45752811.060000002

# Row 16 - D only
$ws.Range("D16").Value = 50601311.960000001

# Row 18 totals - C18 becomes a formula (was hard-coded), D18 stays a formula
$ws.Range("C18").Formula = "=SUM(C12:C17)"

# Row 19
$ws.Range("C19").Value = -392700000.00000024
$ws.Range("D19").Value = 383100000

# Row 21 totals - C21 becomes a formula (was hard-coded), D21 stays a formula
$ws.Range("C21").Formula = "=SUM(C18:C20)"

# Row 22 - D only
$ws.Range("D22").Value = 32201025

# Row 26 - base values
$ws.Range("C26").Value = 1003368420.8613656
$ws.Range("D26").Value = 1009991810

$excel.CalculateFull()
